$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 22, pushing existing rows 22:31 down to 23:32
$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22 with the new weekly record
$ws.Range("A22").Value = 7
$ws.Range("B22").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C22").Value = "Ñuble"
$ws.Range("D22").Value = "2021-11-10"
$ws.Range("E22").Value = 16
$ws.Range("F22").Value = 100112013
$ws.Range("G22").Value = "Alcachofa"
$ws.Range("H22").Value = "Madrigal"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 120
$ws.Range("K22").Value = 11000
$ws.Range("L22").Value = 12000
$ws.Range("M22").Value = 11500
$ws.Range("N22").Value = "$/caja 40 unidades"
$ws.Range("O22").Value = "Provincia del Elquí"
$ws.Range("P22").Value = 288
$ws.Range("Q22").Value = 40
$ws.Range("R22").Value = "Hortaliza"
